$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I62").Value = 11865
$ws.Range("K62").Value = 11865
$ws.Range("H62").Value = 5038.222
$ws.Range("M62").Value = -11241
$ws.Range("I65").Value = 11865
$ws.Range("K65").Value = 59325
$ws.Range("H65").Value = 5038.222
$ws.Range("M65").Value = -56205
$ws.Range("M111").Value = -1433
$ws.Range("J111").Value = 1500
$ws.Range("I111").Value = 1500
$ws.Range("K111").Value = 4500
$ws.Range("H111").Value = 1500
$ws.Range("N111").Value = -10634
$ws.Range("L111").Value = 4500
$ws.Range("J129").Value = 1972.9678
$ws.Range("I129").Value = 5325.6665
$ws.Range("K129").Value = 15976.9995
$ws.Range("H129").Value = 2268.7942
$ws.Range("M129").Value = -10976.9995
$ws.Range("L129").Value = 5918.903399999999
$ws.Range("N129").Value = -15918.9034
$ws.Range("I138").Value = 940.12823
$ws.Range("K138").Value = 2820.38469
$ws.Range("H138").Value = 2464.22
$ws.Range("N138").Value = -20595.9182
$ws.Range("M138").Value = 2319.61531
$ws.Range("L138").Value = 10315.9182
$ws.Range("J138").Value = 3438.6394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 9512.712
$ws.Range("K32").Value = 9512.712
$ws.Range("H32").Value = 10145.84
$ws.Range("M32").Value = -9225.712
$ws.Range("I41").Value = 3800
$ws.Range("K41").Value = 3800
$ws.Range("H41").Value = 3800
$ws.Range("M41").Value = -3386
$ws.Range("M45").Value = -1294.2222
$ws.Range("I45").Value = 1671.2222
$ws.Range("K45").Value = 1671.2222
$ws.Range("H45").Value = 2235.8333
$ws.Range("M122").Value = -1864.6819
$ws.Range("L122").Value = 3684.1875
$ws.Range("I122").Value = 1438.2273
$ws.Range("K122").Value = 4314.6819
$ws.Range("J122").Value = 1228.0625
$ws.Range("H122").Value = 1382.1833
$ws.Range("N122").Value = -8584.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J134").Value = 3579.0527
$ws.Range("I134").Value = 914.65216
$ws.Range("K134").Value = 2743.95648
$ws.Range("H134").Value = 2389.1262
$ws.Range("N134").Value = -15807.1581
$ws.Range("L134").Value = 10737.1581
$ws.Range("M134").Value = -208.9564799999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 889.6667
$ws.Range("H16").Value = 1027.7142
$ws.Range("N16").Value = -1463.6667
$ws.Range("I16").Value = 1131.25
$ws.Range("K16").Value = 1131.25
$ws.Range("M16").Value = -844.25
$ws.Range("L16").Value = 889.6667
$ws.Range("H31").Value = 165920.27
$ws.Range("N31").Value = -166510.27
$ws.Range("L31").Value = 165920.27
$ws.Range("M31").ClearContents()
$ws.Range("J31").Value = 165920.27
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("H34").Value = 165920.27
$ws.Range("N34").Value = -166324.27
$ws.Range("M34").ClearContents()
$ws.Range("L34").Value = 165920.27
$ws.Range("J34").Value = 165920.27
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("J113").Value = 889.6667
$ws.Range("I113").Value = 1131.25
$ws.Range("K113").Value = 1131.25
$ws.Range("H113").Value = 1027.7142
$ws.Range("N113").Value = -5229.6667
$ws.Range("M113").Value = 1038.75
$ws.Range("L113").Value = 889.6667
$ws.Range("J132").Value = 159558.33
$ws.Range("I132").Value = 1321.7073
$ws.Range("K132").Value = 3965.1219
$ws.Range("H132").Value = 29804.3
$ws.Range("N132").Value = -483734.99
$ws.Range("M132").Value = -1435.1219
$ws.Range("L132").Value = 478674.99

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J58").Value = 1458754.4
$ws.Range("H58").Value = 1276532.6
$ws.Range("N58").Value = -4376519.199999999
$ws.Range("L58").Value = 4376263.199999999
$ws.Range("J112").Value = 3270.9092
$ws.Range("N112").Value = -12028.7276
$ws.Range("I112").Value = 1666.3334
$ws.Range("K112").Value = 4999.0002
$ws.Range("H112").Value = 3168.4893
$ws.Range("M112").Value = -3891.0002
$ws.Range("L112").Value = 9812.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L33").Value = 22679.334
$ws.Range("J33").Value = 22679.334
$ws.Range("H33").Value = 22679.334
$ws.Range("N33").Value = -23183.334
$ws.Range("H36").Value = 5000
$ws.Range("N36").Value = -5970
$ws.Range("L36").Value = 5000
$ws.Range("J36").Value = 5000
$ws.Range("L64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("H64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("L67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("M102").Value = 464.3846000000001
$ws.Range("L102").Value = 1200
$ws.Range("J102").Value = 1200
$ws.Range("I102").Value = 1157.6154
$ws.Range("K102").Value = 1157.6154
$ws.Range("H102").Value = 1160.6428
$ws.Range("N102").Value = -4444
$ws.Range("J113").Value = 1475.1428
$ws.Range("I113").Value = 1282.5834
$ws.Range("K113").Value = 1282.5834
$ws.Range("H113").Value = 1353.5264
$ws.Range("N113").Value = -5815.1428
$ws.Range("M113").Value = 887.4166
$ws.Range("L113").Value = 1475.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J40").Value = 3925.2856
$ws.Range("I40").Value = 1631.625
$ws.Range("K40").Value = 1631.625
$ws.Range("H40").Value = 2149.5483
$ws.Range("N40").Value = -4197.2856
$ws.Range("L40").Value = 3925.2856
$ws.Range("M40").Value = -1495.625
$ws.Range("I61").Value = 5395
$ws.Range("K61").Value = 5395
$ws.Range("H61").Value = 3536.6667
$ws.Range("N61").Value = -3011.5
$ws.Range("M61").Value = -5193
$ws.Range("L61").Value = 2607.5
$ws.Range("J61").Value = 2607.5
$ws.Range("J113").Value = 2607.5
$ws.Range("I113").Value = 5395
$ws.Range("K113").Value = 5395
$ws.Range("H113").Value = 3536.6667
$ws.Range("N113").Value = -6947.5
$ws.Range("M113").Value = -3225
$ws.Range("L113").Value = 2607.5
$ws.Range("M122").Value = -252878.992
$ws.Range("L122").Value = 7330.7145
$ws.Range("I122").Value = 85109.664
$ws.Range("K122").Value = 255328.992
$ws.Range("J122").Value = 2443.5715
$ws.Range("H122").Value = 54653.74
$ws.Range("N122").Value = -12230.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 1650
$ws.Range("K62").Value = 1650
$ws.Range("H62").Value = 3328.5715
$ws.Range("M62").Value = -1026
$ws.Range("I65").Value = 1650
$ws.Range("K65").Value = 8250
$ws.Range("H65").Value = 3328.5715
$ws.Range("M65").Value = -5130
$ws.Range("M100").Value = -432.7143
$ws.Range("L100").Value = 1100
$ws.Range("J100").Value = 550
$ws.Range("I100").Value = 486.85715
$ws.Range("K100").Value = 973.7143
$ws.Range("H100").Value = 500.8889
$ws.Range("N100").Value = -2182
$ws.Range("J113").Value = 660.4286
$ws.Range("I113").Value = 567.2778
$ws.Range("K113").Value = 1701.8334
$ws.Range("H113").Value = 593.36
$ws.Range("N113").Value = -6321.2858
$ws.Range("M113").Value = 468.1666
$ws.Range("L113").Value = 1981.2858
$ws.Range("J132").Value = 2064.7058
$ws.Range("I132").Value = 1865.1794
$ws.Range("K132").Value = 5595.5382
$ws.Range("H132").Value = 1925.75
$ws.Range("N132").Value = -11254.1174
$ws.Range("M132").Value = -3065.5382
$ws.Range("L132").Value = 6194.117400000001
